# Carga_Pasaje_AF.xlsx — "Add files via upload" edit
#
# The source sheet "C_11" is renamed to "C_19". Renaming the sheet through
# the COM object model automatically keeps the workbook's defined name
# (_xlnm._FilterDatabase, which is scoped to this sheet and refers to
# C_11!$B$4:$G$4) in sync, so it becomes C_19!$B$4:$G$4 as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "C_19"

Write-Output "Renamed sheet to: $($ws.Name)"
